# Fruta / hortaliza, semanal
# Update the weekly Cilantro price-report values: dates and the associated
# volume/price/origin figures are refreshed for the "Macroferia Regional de
# Talca" market rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 44386
$ws.Range("J2").Value = 200

# --- Row 3 ---
$ws.Range("D3").Value = 44355

# --- Row 4 ---
$ws.Range("D4").Value = 44348

# --- Row 5 ---
$ws.Range("D5").Value = 44362
$ws.Range("K5").Value = 6500
$ws.Range("L5").Value = 6500
$ws.Range("M5").Value = 6500
$ws.Range("N5").Value = "$/caja 36 atados"
$ws.Range("P5").Value = 181
$ws.Range("Q5").Value = 36

# --- Row 6 ---
$ws.Range("D6").Value = 44376
$ws.Range("K6").Value = 6500
$ws.Range("L6").Value = 6500
$ws.Range("M6").Value = 6500
$ws.Range("P6").Value = 181

# --- Row 8 ---
$ws.Range("D8").Value = 44369
$ws.Range("J8").Value = 100
$ws.Range("N8").Value = "$/caja 20 docenas"
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 7000
$ws.Range("Q8").Value = 1

# --- Row 9 ---
$ws.Range("D9").Value = 44340

# --- Row 10 ---
$ws.Range("D10").Value = 44371

# --- Row 11 ---
$ws.Range("D11").Value = 44342
$ws.Range("J11").Value = 150
$ws.Range("O11").Value = "Región del Maule"

# --- Row 12 ---
$ws.Range("D12").Value = 44364
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 7000
$ws.Range("P12").Value = 194

# --- Row 13 ---
$ws.Range("D13").Value = 44372
$ws.Range("O13").Value = "Región Metropolitana"

# --- Row 15 ---
$ws.Range("D15").Value = 44354
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 194
